$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the crypto price/volume/hour data to reflect the latest scrape.
# Numeric-looking values are stored as text in this sheet (same as the
# source data), so a leading apostrophe is used to force text entry and
# avoid Excel auto-converting them to numbers/dates.

$ws.Range("D2").Value = "'245.64"
$ws.Range("G2").Value = "'13"
$ws.Range("G3").Value = "'13"
$ws.Range("D4").Value = "'5.327"
$ws.Range("G4").Value = "'13"
$ws.Range("G5").Value = "'13"
$ws.Range("D6").Value = "'6.501"
$ws.Range("G6").Value = "'13"
$ws.Range("D7").Value = "'3.140"
$ws.Range("G7").Value = "'13"
$ws.Range("D8").Value = "'0.8170"
$ws.Range("G8").Value = "'13"
$ws.Range("D9").Value = "'0.8671"
$ws.Range("G9").Value = "'13"
$ws.Range("B10").Value = "'One"
$ws.Range("C10").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.01012"
$ws.Range("E10").Value = "'9OneONEBestin24h"
$ws.Range("G10").Value = "'13"
$ws.Range("B11").Value = "'WazirX"
$ws.Range("C11").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1369"
$ws.Range("E11").Value = "'10WazirXWRX"
$ws.Range("G11").Value = "'13"
$ws.Range("B12").Value = "'MandalaExchangeToken"
$ws.Range("C12").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07000"
$ws.Range("E12").Value = "'11MandalaExchangeTokenMDX"
$ws.Range("G12").Value = "'13"
$ws.Range("B13").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "'0.03235"
$ws.Range("E13").Value = "'12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G13").Value = "'13"
$ws.Range("B14").Value = "'BitrueCoin"
$ws.Range("C14").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.02896"
$ws.Range("E14").Value = "'13BitrueCoinBTR"
$ws.Range("G14").Value = "'13"
$ws.Range("B15").Value = "'BitMartToken"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09375"
$ws.Range("E15").Value = "'14BitMartTokenBMX"
$ws.Range("G15").Value = "'13"
$ws.Range("B16").Value = "'MCDex"
$ws.Range("C16").Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "'3.745"
$ws.Range("E16").Value = "'15MCDexMCB"
$ws.Range("G16").Value = "'13"
$ws.Range("B17").Value = "'BitForexToken"
$ws.Range("C17").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "'0.001525"
$ws.Range("E17").Value = "'16BitForexTokenBF"
$ws.Range("G17").Value = "'13"
$ws.Range("B18").Value = "'CoinExToken"
$ws.Range("C18").Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "'0.04700"
$ws.Range("E18").Value = "'17CoinExTokenCET"
$ws.Range("G18").Value = "'13"
$ws.Range("D19").Value = "'0.006185"
$ws.Range("G19").Value = "'13"
$ws.Range("D20").Value = "'0.001239"
$ws.Range("G20").Value = "'13"
$ws.Range("D21").Value = "'0.003841"
$ws.Range("G21").Value = "'13"
$ws.Range("D22").Value = "'0.00008794"
$ws.Range("G22").Value = "'13"
$ws.Range("D23").Value = "'3.543"
$ws.Range("G23").Value = "'13"
$ws.Range("D24").Value = "'2.150"
$ws.Range("G24").Value = "'13"
$ws.Range("D25").Value = "'0.3173"
$ws.Range("G25").Value = "'13"
$ws.Range("G26").Value = "'13"
$ws.Range("D27").Value = "'0.1327"
$ws.Range("G27").Value = "'13"
$ws.Range("D28").Value = "'0.0003012"
$ws.Range("E28").Value = "'27UpBotsUBXT"
$ws.Range("G28").Value = "'13"
$ws.Range("G29").Value = "'13"
$ws.Range("G30").Value = "'13"
$ws.Range("G31").Value = "'13"
$ws.Range("G32").Value = "'13"
$ws.Range("G33").Value = "'13"
$ws.Range("G34").Value = "'13"
$ws.Range("G35").Value = "'13"
$ws.Range("G36").Value = "'13"
$ws.Range("G37").Value = "'13"
$ws.Range("G38").Value = "'13"
$ws.Range("G39").Value = "'13"
$ws.Range("G40").Value = "'13"
$ws.Range("D41").Value = "'0.006426"
$ws.Range("G41").Value = "'13"
$ws.Range("D42").Value = "'0.1054"
$ws.Range("G42").Value = "'13"
$ws.Range("G43").Value = "'13"
$ws.Range("D44").Value = "'0.008678"
$ws.Range("G44").Value = "'13"
$ws.Range("D45").Value = "'0.00005353"
$ws.Range("G45").Value = "'13"
$ws.Range("G46").Value = "'13"
$ws.Range("D47").Value = "'0.3882"
$ws.Range("G47").Value = "'13"
$ws.Range("D48").Value = "'0.002569"
$ws.Range("E48").Value = "'47BOLOBOLOWorstin24h"
$ws.Range("G48").Value = "'13"
$ws.Range("G49").Value = "'13"
$ws.Range("G50").Value = "'13"
$ws.Range("G51").Value = "'13"

Write-Host "Updated symbol list"
